$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append "ELK0003" as a new line to the ELK code list cell (merged B14:B17)
$b14 = $ws.Range("B14")
$b14.Value2 = $b14.Value2 + "`n" + "ELK0003"

# Append "Node & Cct Deletion (DN)" as a new line to E17
$e17 = $ws.Range("E17")
$e17.Value2 = $e17.Value2 + "`n" + "Node & Cct Deletion (DN)"
$e17.WrapText = $true
